# Update "F" column (想去人数 / interest count) values across three sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 633
$ws.Range("F3").Value = 672
$ws.Range("F4").Value = 931
$ws.Range("F5").Value = 695
$ws.Range("F6").Value = 828
$ws.Range("F7").Value = 387
$ws.Range("F8").Value = 586
$ws.Range("F10").Value = 1189
$ws.Range("F11").Value = 619
$ws.Range("F12").Value = 368
$ws.Range("F13").Value = 491
$ws.Range("F14").Value = 163
$ws.Range("F15").Value = 318
$ws.Range("F16").Value = 329
$ws.Range("F19").Value = 545
$ws.Range("F20").Value = 60
$ws.Range("F21").Value = 559
$ws.Range("F23").Value = 663

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 78
$ws.Range("F3").Value = 59
$ws.Range("F4").Value = 310
$ws.Range("F9").Value = 218
$ws.Range("F10").Value = 47
$ws.Range("F11").Value = 22
$ws.Range("F13").Value = 58

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 78
$ws.Range("F4").Value = 633
$ws.Range("F5").Value = 59
$ws.Range("F6").Value = 310
$ws.Range("F7").Value = 672
$ws.Range("F8").Value = 931
$ws.Range("F9").Value = 695
$ws.Range("F10").Value = 828
$ws.Range("F11").Value = 387
$ws.Range("F12").Value = 586
$ws.Range("F14").Value = 1189
$ws.Range("F15").Value = 619
$ws.Range("F18").Value = 368
$ws.Range("F19").Value = 491
$ws.Range("F21").Value = 163
$ws.Range("F22").Value = 318
$ws.Range("F24").Value = 329
$ws.Range("F27").Value = 218
$ws.Range("F28").Value = 47
$ws.Range("F29").Value = 545
$ws.Range("F30").Value = 22
$ws.Range("F32").Value = 58
$ws.Range("F33").Value = 60
$ws.Range("F34").Value = 559
$ws.Range("F36").Value = 663
